$d = $word.ActiveDocument

$baseCount = $d.Paragraphs.Count
$insertPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)

$xml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="linespace"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="linespace"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="linespace"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:tab/></w:r><w:r><w:t xml:space="preserve">IVAN </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>McKEE</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:tab/><w:t>Authorised to sign by the Scottish Ministers</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SigBlock"/><w:rPr><w:rStyle w:val="Sigtitle"/></w:rPr></w:pPr><w:r><w:t>St Andrew’s House,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:t>Edinburgh</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:t>1st October 2024</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="linespace"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:t>We consent</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="linespace"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="SigBlock"/><w:rPr><w:rStyle w:val="SigSignee"/></w:rPr></w:pPr><w:r><w:tab/></w:r><w:r><w:t>JEFF SMITH</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:tab/></w:r><w:r><w:t>ANNA TURLEY</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SigBlock"/><w:rPr><w:rStyle w:val="Sigtitle"/></w:rPr></w:pPr><w:r><w:tab/><w:t>Two of the Lords Commissioners of His Majesty’s Treasury</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="LQN2"/><w:ind w:left="0" w:firstLine="0"/></w:pPr><w:r><w:t>1st October 2024</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertPoint.InsertXML($xml)

# InsertXML faithfully restores paragraph styles, tabs, proofErr markers, and
# even rStyle references that live in a paragraph's pPr/rPr (paragraph-mark
# formatting) -- but it silently drops rStyle references that live on the
# w:rPr of individual w:r runs. Re-apply those character styles now, by
# addressing each new paragraph's Range (relative indices, so this is robust
# to $baseCount).

function Get-NewParagraph([int]$offset) {
    return $d.Paragraphs.Item($baseCount + $offset)
}

# 4: tab + "IVAN " + "McKEE"  -> style "IVAN McKEE" (i.e. skip the leading tab) as SigSignee
$p = Get-NewParagraph 4
$full = $p.Range
$sub = $d.Range($full.Start + 1, $full.End - 1)
$sub.Style = "SigSignee"

# 5: tab + "Authorised to sign by the Scottish Ministers" (one run) -> Sigtitle, include the tab
$p = Get-NewParagraph 5
$full = $p.Range
$sub = $d.Range($full.Start, $full.End - 1)
$sub.Style = "Sigtitle"

# 6: "St Andrew's House," -> SigAdd
$p = Get-NewParagraph 6
$full = $p.Range
$sub = $d.Range($full.Start, $full.End - 1)
$sub.Style = "SigAdd"

# 7: "Edinburgh" -> SigAdd
$p = Get-NewParagraph 7
$full = $p.Range
$sub = $d.Range($full.Start, $full.End - 1)
$sub.Style = "SigAdd"

# 8: "1st October 2024" -> SigDate
$p = Get-NewParagraph 8
$full = $p.Range
$sub = $d.Range($full.Start, $full.End - 1)
$sub.Style = "SigDate"

# 10: "We consent" -> Sigsignatory
$p = Get-NewParagraph 10
$full = $p.Range
$sub = $d.Range($full.Start, $full.End - 1)
$sub.Style = "Sigsignatory"

# 12: tab + "JEFF SMITH" -> style "JEFF SMITH" (skip the leading tab) as SigSignee
$p = Get-NewParagraph 12
$full = $p.Range
$sub = $d.Range($full.Start + 1, $full.End - 1)
$sub.Style = "SigSignee"

# 13: tab + "ANNA TURLEY" -> style "ANNA TURLEY" (skip the leading tab) as SigSignee
$p = Get-NewParagraph 13
$full = $p.Range
$sub = $d.Range($full.Start + 1, $full.End - 1)
$sub.Style = "SigSignee"

# 14: tab + "Two of the Lords Commissioners of His Majesty's Treasury" (one run) -> Sigtitle, include the tab
$p = Get-NewParagraph 14
$full = $p.Range
$sub = $d.Range($full.Start, $full.End - 1)
$sub.Style = "Sigtitle"

# 15: "1st October 2024" (LQN2 paragraph) -> SigDate
$p = Get-NewParagraph 15
$full = $p.Range
$sub = $d.Range($full.Start, $full.End - 1)
$sub.Style = "SigDate"

Write-Host "Paragraphs after edit: $($d.Paragraphs.Count)"
